$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value  = 4
$ws.Range("E4").Value  = 5
$ws.Range("E5").Value  = 0
$ws.Range("E6").Value  = 1
$ws.Range("E7").Value  = 5
$ws.Range("E8").Value  = 5
$ws.Range("E9").Value  = 4
$ws.Range("E10").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("E13").Value = 5
$ws.Range("E15").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("E21").Value = 4
